$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 178.14285  # H33: 147.6 -> 178.14285
$ws.Cells.Item(33, 9).Value = 178.14285  # I33: 156.77777 -> 178.14285
$ws.Cells.Item(33, 10).Value = 0  # J33: 65 -> 0
$ws.Cells.Item(33, 11).Value = 178.14285  # K33: 156.77777 -> 178.14285
$ws.Cells.Item(33, 12).Value = 0  # L33: 65 -> 0
$ws.Cells.Item(33, 13).Value = 50.85714999999999  # M33: 72.22223 -> 50.85714999999999
$ws.Cells.Item(33, 14).ClearContents()  # N33: delete (was -523)
$ws.Cells.Item(38, 8).Value = 5814.6665  # H38: 4082.6191 -> 5814.6665
$ws.Cells.Item(38, 9).Value = 2430.375  # I38: 1645.25 -> 2430.375
$ws.Cells.Item(38, 10).Value = 8522.1  # J38: 7332.4443 -> 8522.1
$ws.Cells.Item(38, 11).Value = 7291.125  # K38: 4935.75 -> 7291.125
$ws.Cells.Item(38, 12).Value = 25566.3  # L38: 21997.3329 -> 25566.3
$ws.Cells.Item(38, 13).Value = -6919.125  # M38: -4563.75 -> -6919.125
$ws.Cells.Item(38, 14).Value = -26310.3  # N38: -22741.3329 -> -26310.3
$ws.Cells.Item(40, 8).Value = 11807  # H40: 13260.714 -> 11807
$ws.Cells.Item(40, 9).Value = 9137.5  # I40: 10057.143 -> 9137.5
$ws.Cells.Item(40, 10).Value = 16552.777  # J40: 19667.857 -> 16552.777
$ws.Cells.Item(40, 11).Value = 9137.5  # K40: 10057.143 -> 9137.5
$ws.Cells.Item(40, 12).Value = 16552.777  # L40: 19667.857 -> 16552.777
$ws.Cells.Item(40, 13).Value = -8962.5  # M40: -9882.143 -> -8962.5
$ws.Cells.Item(40, 14).Value = -16902.777  # N40: -20017.857 -> -16902.777
$ws.Cells.Item(61, 8).Value = 341.5  # H61: 338.66666 -> 341.5
$ws.Cells.Item(61, 9).Value = 341.5  # I61: 338.66666 -> 341.5
$ws.Cells.Item(61, 11).Value = 1024.5  # K61: 1015.99998 -> 1024.5
$ws.Cells.Item(61, 13).Value = -852.5  # M61: -843.9999799999999 -> -852.5
$ws.Cells.Item(113, 8).Value = 3569.7144  # H113: 3594.5557 -> 3569.7144
$ws.Cells.Item(113, 9).Value = 3096  # I113: 3109.1333 -> 3096
$ws.Cells.Item(113, 11).Value = 3096  # K113: 3109.1333 -> 3096
$ws.Cells.Item(113, 13).Value = 158  # M113: 144.8667 -> 158
$ws.Cells.Item(121, 8).Value = 4764.5386  # H121: 4765.3076 -> 4764.5386
$ws.Cells.Item(121, 10).Value = 4764.5386  # J121: 4765.3076 -> 4764.5386
$ws.Cells.Item(121, 12).Value = 14293.6158  # L121: 14295.9228 -> 14293.6158
$ws.Cells.Item(121, 14).Value = -17787.6158  # N121: -17789.9228 -> -17787.6158
$ws.Cells.Item(132, 8).Value = 13662.122  # H132: 14824.022 -> 13662.122
$ws.Cells.Item(132, 9).Value = 3099.6667  # I132: 3602.4 -> 3099.6667
$ws.Cells.Item(132, 10).Value = 19795.162  # J132: 20434.834 -> 19795.162
$ws.Cells.Item(132, 11).Value = 9299.000100000001  # K132: 10807.2 -> 9299.000100000001
$ws.Cells.Item(132, 12).Value = 59385.486  # L132: 61304.50199999999 -> 59385.486
$ws.Cells.Item(132, 13).Value = -6769.000100000001  # M132: -8277.200000000001 -> -6769.000100000001
$ws.Cells.Item(132, 14).Value = -64445.486  # N132: -66364.50199999999 -> -64445.486
$ws.Cells.Item(138, 8).Value = 2398.961  # H138: 2250.9893 -> 2398.961
$ws.Cells.Item(138, 9).Value = 1418.0952  # I138: 1362.6818 -> 1418.0952
$ws.Cells.Item(138, 10).Value = 2766.7856  # J138: 2522.4167 -> 2766.7856
$ws.Cells.Item(138, 11).Value = 4254.2856  # K138: 4088.0454 -> 4254.2856
$ws.Cells.Item(138, 12).Value = 8300.356800000001  # L138: 7567.250100000001 -> 8300.356800000001
$ws.Cells.Item(138, 13).Value = 885.7143999999998  # M138: 1051.9546 -> 885.7143999999998
$ws.Cells.Item(138, 14).Value = -18580.3568  # N138: -17847.2501 -> -18580.3568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 21295.09  # H32: 17815.092 -> 21295.09
$ws.Cells.Item(32, 9).Value = 23890.969  # I32: 19798.225 -> 23890.969
$ws.Cells.Item(32, 10).Value = 14156.417  # J32: 12149 -> 14156.417
$ws.Cells.Item(32, 11).Value = 23890.969  # K32: 19798.225 -> 23890.969
$ws.Cells.Item(32, 12).Value = 14156.417  # L32: 12149 -> 14156.417
$ws.Cells.Item(32, 13).Value = -23603.969  # M32: -19511.225 -> -23603.969
$ws.Cells.Item(32, 14).Value = -14730.417  # N32: -12723 -> -14730.417
$ws.Cells.Item(45, 8).Value = 2957.4167  # H45: 3090 -> 2957.4167
$ws.Cells.Item(45, 9).Value = 1784.8572  # I45: 1832.5 -> 1784.8572
$ws.Cells.Item(45, 11).Value = 1784.8572  # K45: 1832.5 -> 1784.8572
$ws.Cells.Item(45, 13).Value = -1407.8572  # M45: -1455.5 -> -1407.8572
$ws.Cells.Item(61, 8).Value = 9829.956  # H61: 8611.654 -> 9829.956
$ws.Cells.Item(61, 9).Value = 10449.389  # I61: 8722.862999999999 -> 10449.389
$ws.Cells.Item(61, 10).Value = 7600  # J61: 8000 -> 7600
$ws.Cells.Item(61, 11).Value = 10449.389  # K61: 8722.862999999999 -> 10449.389
$ws.Cells.Item(61, 12).Value = 7600  # L61: 8000 -> 7600
$ws.Cells.Item(61, 13).Value = -10237.389  # M61: -8510.862999999999 -> -10237.389
$ws.Cells.Item(61, 14).Value = -8024  # N61: -8424 -> -8024
$ws.Cells.Item(63, 8).Value = 3832.3333  # H63: 4018.8 -> 3832.3333
$ws.Cells.Item(63, 10).Value = 4723.5  # J63: 5331.3335 -> 4723.5
$ws.Cells.Item(63, 12).Value = 4723.5  # L63: 5331.3335 -> 4723.5
$ws.Cells.Item(63, 14).Value = -6095.5  # N63: -6703.3335 -> -6095.5
$ws.Cells.Item(66, 8).Value = 3832.3333  # H66: 4018.8 -> 3832.3333
$ws.Cells.Item(66, 10).Value = 4723.5  # J66: 5331.3335 -> 4723.5
$ws.Cells.Item(66, 12).Value = 23617.5  # L66: 26656.6675 -> 23617.5
$ws.Cells.Item(66, 14).Value = -30481.5  # N66: -33520.6675 -> -30481.5
$ws.Cells.Item(97, 8).Value = 315.89285  # H97: 311.82758 -> 315.89285
$ws.Cells.Item(97, 9).Value = 328.61905  # I97: 322.68182 -> 328.61905
$ws.Cells.Item(97, 11).Value = 328.61905  # K97: 322.68182 -> 328.61905
$ws.Cells.Item(97, 13).Value = 167.38095  # M97: 173.31818 -> 167.38095
$ws.Cells.Item(102, 8).Value = 490810.16  # H102: 528507.0600000001 -> 490810.16
$ws.Cells.Item(102, 9).Value = 623799.25  # I102: 653485.4399999999 -> 623799.25
$ws.Cells.Item(102, 10).Value = 3183.3333  # J102: 3597.8 -> 3183.3333
$ws.Cells.Item(102, 11).Value = 623799.25  # K102: 653485.4399999999 -> 623799.25
$ws.Cells.Item(102, 12).Value = 3183.3333  # L102: 3597.8 -> 3183.3333
$ws.Cells.Item(102, 13).Value = -622177.25  # M102: -651863.4399999999 -> -622177.25
$ws.Cells.Item(102, 14).Value = -6427.3333  # N102: -6841.8 -> -6427.3333
$ws.Cells.Item(122, 8).Value = 4832.4644  # H122: 4959.926 -> 4832.4644
$ws.Cells.Item(122, 9).Value = 2691  # I122: 2767.4707 -> 2691
$ws.Cells.Item(122, 11).Value = 8073  # K122: 8302.4121 -> 8073
$ws.Cells.Item(122, 13).Value = -5623  # M122: -5852.4121 -> -5623
$ws.Cells.Item(132, 8).Value = 12085.884  # H132: 12683.368 -> 12085.884
$ws.Cells.Item(132, 9).Value = 15971.897  # I132: 17241.75 -> 15971.897
$ws.Cells.Item(132, 11).Value = 47915.69100000001  # K132: 51725.25 -> 47915.69100000001
$ws.Cells.Item(132, 13).Value = -45385.69100000001  # M132: -49195.25 -> -45385.69100000001
$ws.Cells.Item(136, 8).Value = 9829.956  # H136: 8611.654 -> 9829.956
$ws.Cells.Item(136, 9).Value = 10449.389  # I136: 8722.862999999999 -> 10449.389
$ws.Cells.Item(136, 10).Value = 7600  # J136: 8000 -> 7600
$ws.Cells.Item(136, 11).Value = 31348.167  # K136: 26168.589 -> 31348.167
$ws.Cells.Item(136, 12).Value = 22800  # L136: 24000 -> 22800
$ws.Cells.Item(136, 13).Value = -28798.167  # M136: -23618.589 -> -28798.167
$ws.Cells.Item(136, 14).Value = -27900  # N136: -29100 -> -27900

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 30999  # H7: 49999 -> 30999
$ws.Cells.Item(7, 9).Value = 11999  # I7: 0 -> 11999
$ws.Cells.Item(7, 11).Value = 11999  # K7: 0 -> 11999
$ws.Cells.Item(7, 13).Value = -11886  # M7: None -> -11886
$ws.Cells.Item(82, 8).Value = 50000  # H82: 17628.5 -> 50000
$ws.Cells.Item(82, 9).Value = 50000  # I82: 17628.5 -> 50000
$ws.Cells.Item(82, 11).Value = 50000  # K82: 17628.5 -> 50000
$ws.Cells.Item(82, 13).Value = -49617  # M82: -17245.5 -> -49617
$ws.Cells.Item(85, 8).Value = 50000  # H85: 17628.5 -> 50000
$ws.Cells.Item(85, 9).Value = 50000  # I85: 17628.5 -> 50000
$ws.Cells.Item(85, 11).Value = 50000  # K85: 17628.5 -> 50000
$ws.Cells.Item(85, 13).Value = -48674  # M85: -16302.5 -> -48674
$ws.Cells.Item(99, 8).Value = 1158794.5  # H99: 1303545.8 -> 1158794.5
$ws.Cells.Item(99, 9).Value = 1489392.9  # I99: 1737494.2 -> 1489392.9
$ws.Cells.Item(99, 11).Value = 1489392.9  # K99: 1737494.2 -> 1489392.9
$ws.Cells.Item(99, 13).Value = -1487894.9  # M99: -1735996.2 -> -1487894.9
$ws.Cells.Item(105, 8).Value = 2470.4707  # H105: 2499.9375 -> 2470.4707
$ws.Cells.Item(105, 9).Value = 2437.4375  # I105: 2466.6667 -> 2437.4375
$ws.Cells.Item(105, 11).Value = 2437.4375  # K105: 2466.6667 -> 2437.4375
$ws.Cells.Item(105, 13).Value = -690.4375  # M105: -719.6667000000002 -> -690.4375
$ws.Cells.Item(134, 8).Value = 1162.25  # H134: 1232.6364 -> 1162.25
$ws.Cells.Item(134, 9).Value = 1175.1  # I134: 1215.9474 -> 1175.1
$ws.Cells.Item(134, 10).Value = 1098  # J134: 1338.3334 -> 1098
$ws.Cells.Item(134, 11).Value = 3525.3  # K134: 3647.8422 -> 3525.3
$ws.Cells.Item(134, 12).Value = 3294  # L134: 4015.0002 -> 3294
$ws.Cells.Item(134, 13).Value = -990.2999999999997  # M134: -1112.8422 -> -990.2999999999997
$ws.Cells.Item(134, 14).Value = -8364  # N134: -9085.0002 -> -8364

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(93, 8).Value = 32000  # H93: 250005500 -> 32000
$ws.Cells.Item(93, 9).Value = 5000  # I93: 7331.3335 -> 5000
$ws.Cells.Item(93, 10).Value = 59000  # J93: 1000000000 -> 59000
$ws.Cells.Item(93, 11).Value = 5000  # K93: 7331.3335 -> 5000
$ws.Cells.Item(93, 12).Value = 59000  # L93: 1000000000 -> 59000
$ws.Cells.Item(93, 13).Value = -3128  # M93: -5459.3335 -> -3128
$ws.Cells.Item(93, 14).Value = -62744  # N93: -1000003744 -> -62744
$ws.Cells.Item(103, 8).Value = 28997.8  # H103: 47333 -> 28997.8
$ws.Cells.Item(103, 9).Value = 8378.5  # I103: 15262 -> 8378.5
$ws.Cells.Item(103, 11).Value = 8378.5  # K103: 15262 -> 8378.5
$ws.Cells.Item(103, 13).Value = -7206.5  # M103: -14090 -> -7206.5
$ws.Cells.Item(132, 8).Value = 17558318  # H132: 8779635 -> 17558318
$ws.Cells.Item(132, 9).Value = 20850108  # I132: 9531995 -> 20850108
$ws.Cells.Item(132, 11).Value = 62550324  # K132: 28595985 -> 62550324
$ws.Cells.Item(132, 13).Value = -62547794  # M132: -28593455 -> -62547794
$ws.Cells.Item(134, 8).Value = 2562  # H134: 2288.8572 -> 2562
$ws.Cells.Item(134, 10).Value = 0  # J134: 650 -> 0
$ws.Cells.Item(134, 12).Value = 0  # L134: 1950 -> 0
$ws.Cells.Item(134, 14).ClearContents()  # N134: delete (was -7020)

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 5375  # H80: 5714.8 -> 5375
$ws.Cells.Item(80, 9).Value = 3500  # I80: 0 -> 3500
$ws.Cells.Item(80, 10).Value = 6000  # J80: 5714.8 -> 6000
$ws.Cells.Item(80, 11).Value = 10500  # K80: 0 -> 10500
$ws.Cells.Item(80, 12).Value = 18000  # L80: 17144.4 -> 18000
$ws.Cells.Item(80, 13).Value = -9564  # M80: None -> -9564
$ws.Cells.Item(80, 14).Value = -19872  # N80: -19016.4 -> -19872
$ws.Cells.Item(83, 8).Value = 5375  # H83: 5714.8 -> 5375
$ws.Cells.Item(83, 9).Value = 3500  # I83: 0 -> 3500
$ws.Cells.Item(83, 10).Value = 6000  # J83: 5714.8 -> 6000
$ws.Cells.Item(83, 11).Value = 31500  # K83: 0 -> 31500
$ws.Cells.Item(83, 12).Value = 54000  # L83: 51433.2 -> 54000
$ws.Cells.Item(83, 13).Value = -26820  # M83: None -> -26820
$ws.Cells.Item(83, 14).Value = -63360  # N83: -60793.2 -> -63360
$ws.Cells.Item(138, 8).Value = 4412  # H138: 5093.727 -> 4412
$ws.Cells.Item(138, 9).Value = 6500  # I138: 10555.5 -> 6500
$ws.Cells.Item(138, 10).Value = 3890  # J138: 3880 -> 3890
$ws.Cells.Item(138, 11).Value = 19500  # K138: 31666.5 -> 19500
$ws.Cells.Item(138, 12).Value = 11670  # L138: 11640 -> 11670
$ws.Cells.Item(138, 13).Value = -14360  # M138: -26526.5 -> -14360
$ws.Cells.Item(138, 14).Value = -21950  # N138: -21920 -> -21950

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 13749.923  # H102: 31262156 -> 13749.923
$ws.Cells.Item(102, 9).Value = 16218.75  # I102: 50013350 -> 16218.75
$ws.Cells.Item(102, 10).Value = 9799.799999999999  # J102: 10166.5 -> 9799.799999999999
$ws.Cells.Item(102, 11).Value = 16218.75  # K102: 50013350 -> 16218.75
$ws.Cells.Item(102, 12).Value = 9799.799999999999  # L102: 10166.5 -> 9799.799999999999
$ws.Cells.Item(102, 13).Value = -14596.75  # M102: -50011728 -> -14596.75
$ws.Cells.Item(102, 14).Value = -13043.8  # N102: -13410.5 -> -13043.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 2000  # H2: 10000 -> 2000
$ws.Cells.Item(2, 9).Value = 1000  # I2: 0 -> 1000
$ws.Cells.Item(2, 10).Value = 2666.6667  # J2: 10000 -> 2666.6667
$ws.Cells.Item(2, 11).Value = 1000  # K2: 0 -> 1000
$ws.Cells.Item(2, 12).Value = 2666.6667  # L2: 10000 -> 2666.6667
$ws.Cells.Item(2, 13).Value = -888  # M2: None -> -888
$ws.Cells.Item(2, 14).Value = -2890.6667  # N2: -10224 -> -2890.6667
$ws.Cells.Item(46, 8).Value = 5744.3794  # H46: 5582.9 -> 5744.3794
$ws.Cells.Item(46, 9).Value = 4600  # I46: 2750 -> 4600
$ws.Cells.Item(46, 11).Value = 4600  # K46: 2750 -> 4600
$ws.Cells.Item(46, 13).Value = -4412  # M46: -2562 -> -4412
$ws.Cells.Item(55, 8).Value = 340.6154  # H55: 352 -> 340.6154
$ws.Cells.Item(55, 10).Value = 478  # J55: 480.44446 -> 478
$ws.Cells.Item(55, 12).Value = 478  # L55: 480.44446 -> 478
$ws.Cells.Item(55, 14).Value = -824  # N55: -826.4444599999999 -> -824
$ws.Cells.Item(122, 8).Value = 9898.143  # H122: 10529.154 -> 9898.143
$ws.Cells.Item(122, 9).Value = 4091.7778  # I122: 4391.375 -> 4091.7778
$ws.Cells.Item(122, 11).Value = 12275.3334  # K122: 13174.125 -> 12275.3334
$ws.Cells.Item(122, 13).Value = -9825.3334  # M122: -10724.125 -> -9825.3334
$ws.Cells.Item(132, 8).Value = 4070.5454  # H132: 4182.75 -> 4070.5454
$ws.Cells.Item(132, 9).Value = 4026.2856  # I132: 4182.75 -> 4026.2856
$ws.Cells.Item(132, 10).Value = 5000  # J132: 0 -> 5000
$ws.Cells.Item(132, 11).Value = 12078.8568  # K132: 12548.25 -> 12078.8568
$ws.Cells.Item(132, 12).Value = 15000  # L132: 0 -> 15000
$ws.Cells.Item(132, 13).Value = -9548.856800000001  # M132: -10018.25 -> -9548.856800000001
$ws.Cells.Item(132, 14).Value = -20060  # N132: None -> -20060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 7090.067  # H62: 7488.615 -> 7090.067
$ws.Cells.Item(62, 9).Value = 4992.7144  # I62: 5158.1665 -> 4992.7144
$ws.Cells.Item(62, 10).Value = 8925.25  # J62: 9486.143 -> 8925.25
$ws.Cells.Item(62, 11).Value = 4992.7144  # K62: 5158.1665 -> 4992.7144
$ws.Cells.Item(62, 12).Value = 8925.25  # L62: 9486.143 -> 8925.25
$ws.Cells.Item(62, 13).Value = -4368.7144  # M62: -4534.1665 -> -4368.7144
$ws.Cells.Item(62, 14).Value = -10173.25  # N62: -10734.143 -> -10173.25
$ws.Cells.Item(65, 8).Value = 7090.067  # H65: 7488.615 -> 7090.067
$ws.Cells.Item(65, 9).Value = 4992.7144  # I65: 5158.1665 -> 4992.7144
$ws.Cells.Item(65, 10).Value = 8925.25  # J65: 9486.143 -> 8925.25
$ws.Cells.Item(65, 11).Value = 24963.572  # K65: 25790.8325 -> 24963.572
$ws.Cells.Item(65, 12).Value = 44626.25  # L65: 47430.715 -> 44626.25
$ws.Cells.Item(65, 13).Value = -21843.572  # M65: -22670.8325 -> -21843.572
$ws.Cells.Item(65, 14).Value = -50866.25  # N65: -53670.715 -> -50866.25
$ws.Cells.Item(126, 8).Value = 3049.6365  # H126: 2854.8 -> 3049.6365
$ws.Cells.Item(126, 10).Value = 5665.6665  # J126: 5999.5 -> 5665.6665
$ws.Cells.Item(126, 12).Value = 16996.9995  # L126: 17998.5 -> 16996.9995
$ws.Cells.Item(126, 14).Value = -21936.9995  # N126: -22938.5 -> -21936.9995
$ws.Cells.Item(132, 8).Value = 21742782  # H132: 23813514 -> 21742782
$ws.Cells.Item(132, 9).Value = 526.41174  # I132: 584 -> 526.41174
$ws.Cells.Item(132, 11).Value = 1579.23522  # K132: 1752 -> 1579.23522
$ws.Cells.Item(132, 13).Value = 950.76478  # M132: 778 -> 950.76478
